# Updates the crypto price/volume table (and one name swap) to match
# the latest scrape, per the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
${ws}.Range("D2").Value = '27.470.69'
${ws}.Range("E2").Value = '  -3.33%  '

# Row 3
${ws}.Range("D3").Value = '1.755.11'
${ws}.Range("E3").Value = '  -2.72%  '

# Row 4
${ws}.Range("D4").NumberFormat = "@"
${ws}.Range("D4").Value = '1.005'

# Row 5
${ws}.Range("D5").NumberFormat = "@"
${ws}.Range("D5").Value = '323.94'
${ws}.Range("E5").Value = '  -1.32%  '

# Row 6
${ws}.Range("D6").NumberFormat = "@"
${ws}.Range("D6").Value = '1.001'
${ws}.Range("E6").Value = '  +0.17%  '

# Row 7
${ws}.Range("D7").NumberFormat = "@"
${ws}.Range("D7").Value = '0.4393'
${ws}.Range("E7").Value = '  -1.32%  '

# Row 8
${ws}.Range("E8").Value = '  -1.87%  '

# Row 9
${ws}.Range("D9").NumberFormat = "@"
${ws}.Range("D9").Value = '44.82'
${ws}.Range("E9").Value = '  +0.49%  '

# Row 10
${ws}.Range("D10").NumberFormat = "@"
${ws}.Range("D10").Value = '0.07659'
${ws}.Range("E10").Value = '  +2.11%  '

# Row 11
${ws}.Range("D11").NumberFormat = "@"
${ws}.Range("D11").Value = '1.111'
${ws}.Range("E11").Value = '  -3.26%  '

# Row 12
${ws}.Range("D12").NumberFormat = "@"
${ws}.Range("D12").Value = '1.003'
${ws}.Range("E12").Value = '  +0.13%  '

# Row 13
${ws}.Range("D13").NumberFormat = "@"
${ws}.Range("D13").Value = '21.54'
${ws}.Range("E13").Value = '  -4.53%  '

# Row 14
${ws}.Range("D14").NumberFormat = "@"
${ws}.Range("D14").Value = '6.150'
${ws}.Range("E14").Value = '  -2.32%  '

# Row 15
${ws}.Range("D15").NumberFormat = "@"
${ws}.Range("D15").Value = '7.409'
${ws}.Range("E15").Value = '  -2.75%  '

# Row 16
${ws}.Range("D16").Value = '1.763.03'
${ws}.Range("E16").Value = '  -2.18%  '

# Row 17
${ws}.Range("D17").NumberFormat = "@"
${ws}.Range("D17").Value = '90.22'
${ws}.Range("E17").Value = '  +12.02%  '

# Row 18
${ws}.Range("D18").NumberFormat = "@"
${ws}.Range("D18").Value = '0.00001071'
${ws}.Range("E18").Value = '  -1.91%  '

# Row 19
${ws}.Range("D19").NumberFormat = "@"
${ws}.Range("D19").Value = '0.06229'
${ws}.Range("E19").Value = '  -8.43%  '

# Row 20
${ws}.Range("E20").Value = '  +0.09%  '

# Row 21
${ws}.Range("D21").NumberFormat = "@"
${ws}.Range("D21").Value = '17.36'
${ws}.Range("E21").Value = '  -1.37%  '

# Row 22
${ws}.Range("D22").NumberFormat = "@"
${ws}.Range("D22").Value = '6.169'
${ws}.Range("E22").Value = '  -2.49%  '

# Row 23
${ws}.Range("D23").NumberFormat = "@"
${ws}.Range("D23").Value = '0.5274'
${ws}.Range("E23").Value = '  -3.07%  '

# Row 24
${ws}.Range("D24").Value = '27.518.99'
${ws}.Range("E24").Value = '  -3.15%  '

# Row 25
${ws}.Range("D25").NumberFormat = "@"
${ws}.Range("D25").Value = '11.50'
${ws}.Range("E25").Value = '  -2.57%  '

# Row 26
${ws}.Range("D26").NumberFormat = "@"
${ws}.Range("D26").Value = '2.300'
${ws}.Range("E26").Value = '  -4.74%  '

# Row 27
${ws}.Range("D27").NumberFormat = "@"
${ws}.Range("D27").Value = '20.47'
${ws}.Range("E27").Value = '  -0.05%  '

# Row 28
${ws}.Range("D28").NumberFormat = "@"
${ws}.Range("D28").Value = '152.72'
${ws}.Range("E28").Value = '  -0.68%  '

# Row 29
${ws}.Range("D29").NumberFormat = "@"
${ws}.Range("D29").Value = '2.287'
${ws}.Range("E29").Value = '  -2.60%  '

# Row 30
${ws}.Range("D30").Value = '1.956.96'
${ws}.Range("E30").Value = '  -2.47%  '

# Row 31
${ws}.Range("D31").NumberFormat = "@"
${ws}.Range("D31").Value = '127.38'
${ws}.Range("E31").Value = '  -3.72%  '

# Row 32
${ws}.Range("D32").NumberFormat = "@"
${ws}.Range("D32").Value = '1.172'
${ws}.Range("E32").Value = '  -6.61%  '

# Row 33
${ws}.Range("D33").NumberFormat = "@"
${ws}.Range("D33").Value = '5.699'
${ws}.Range("E33").Value = '  -2.00%  '

# Row 34
${ws}.Range("D34").NumberFormat = "@"
${ws}.Range("D34").Value = '0.09162'
${ws}.Range("E34").Value = '  -1.85%  '

# Row 35
${ws}.Range("E35").Value = '  -9.02%  '

# Row 36
${ws}.Range("D36").NumberFormat = "@"
${ws}.Range("D36").Value = '12.55'
${ws}.Range("E36").Value = '  +3.49%  '

# Row 37
${ws}.Range("D37").NumberFormat = "@"
${ws}.Range("D37").Value = '0.02309'
${ws}.Range("E37").Value = '  -1.44%  '

# Row 38
${ws}.Range("D38").NumberFormat = "@"
${ws}.Range("D38").Value = '0.2150'
${ws}.Range("E38").Value = '  -5.40%  '

# Row 39
${ws}.Range("D39").NumberFormat = "@"
${ws}.Range("D39").Value = '0.06098'
${ws}.Range("E39").Value = '  -3.96%  '

# Row 40
${ws}.Range("D40").NumberFormat = "@"
${ws}.Range("D40").Value = '5.046'
${ws}.Range("E40").Value = '  -2.17%  '

# Row 41
${ws}.Range("D41").NumberFormat = "@"
${ws}.Range("D41").Value = '0.6423'
${ws}.Range("E41").Value = '  -2.54%  '

# Row 42
${ws}.Range("D42").NumberFormat = "@"
${ws}.Range("D42").Value = '1.176'
${ws}.Range("E42").Value = '  -2.65%  '

# Row 43
${ws}.Range("B43").Value = 'Frax'
${ws}.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
${ws}.Range("D43").NumberFormat = "@"
${ws}.Range("D43").Value = '1.001'
${ws}.Range("E43").Value = '  +0.20%  '

# Row 44
${ws}.Range("B44").Value = 'FraxShare'
${ws}.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
${ws}.Range("D44").NumberFormat = "@"
${ws}.Range("D44").Value = '7.913'
${ws}.Range("E44").Value = '  -2.17%  '

# Row 45
${ws}.Range("D45").NumberFormat = "@"
${ws}.Range("D45").Value = '1.386'
${ws}.Range("E45").Value = '  -4.52%  '

# Row 46
${ws}.Range("D46").NumberFormat = "@"
${ws}.Range("D46").Value = '13.70'
${ws}.Range("E46").Value = '  -0.83%  '

# Row 47
${ws}.Range("D47").NumberFormat = "@"
${ws}.Range("D47").Value = '0.5951'
${ws}.Range("E47").Value = '  -2.12%  '

# Row 48
${ws}.Range("D48").NumberFormat = "@"
${ws}.Range("D48").Value = '3.716'
${ws}.Range("E48").Value = '  -2.59%  '

# Row 49
${ws}.Range("D49").NumberFormat = "@"
${ws}.Range("D49").Value = '125.88'
${ws}.Range("E49").Value = '  -1.98%  '

# Row 50
${ws}.Range("D50").NumberFormat = "@"
${ws}.Range("D50").Value = '1.975'
${ws}.Range("E50").Value = '  -2.65%  '

# Row 51
${ws}.Range("D51").NumberFormat = "@"
${ws}.Range("D51").Value = '0.06873'
${ws}.Range("E51").Value = '  -3.04%  '
